# Add three new "PUBLONS" test-case rows (9-11) to the "Test Cases" sheet,
# matching the new TCID / JIRA ID / Description / Runmode columns, and move
# the sheet's selection onto the newly populated Runmode column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate() | Out-Null

# Row 9 - PUBLONS006
$ws.Range("A9").Value = "PUBLONS006"
$ws.Range("B9").Value = "OPQA-5779&&OPQA-5780"
$ws.Range("C9").Value = 'Verify that error message "First name is too long." whenever enter more than 50 characters&&Verify that first name should be maximum of 50 characters long and these fields should not be empty.'
$ws.Range("D9").Value = "Y"
$ws.Rows.Item(9).RowHeight = 30

# Row 10 - PUBLONS007
$ws.Range("A10").Value = "PUBLONS007"
$ws.Range("B10").Value = "OPQA-5778"
$ws.Range("C10").Value = 'Verify that error message "Please enter your first name." whenever not enter any text in email field'
$ws.Range("D10").Value = "Y"

# Row 11 - PUBLONS010
$ws.Range("A11").Value = "PUBLONS010"
$ws.Range("B11").Value = "OPQA-5784&&OPQA-5785"
$ws.Range("C11").Value = 'Verify that "Your email address is already registered. Please sign in." error message whenever try to create publons user using existing account.&&Verify that email address prepopulated in sign in page whenever try to register user using existing user'
$ws.Range("D11").Value = "Y"
$ws.Rows.Item(11).RowHeight = 30

# Match the author's final selection: D8:D11 with D8 active.
$ws.Range("D8:D11").Select() | Out-Null
